$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: add a new sample/test-set entry (date, hours, task)
$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B25").Value = 44623

$ws.Range("C25").Value = 2

$ws.Range("D25").Value = "Call e check librerie"

# Row 26: D26 picks up the centered/underlined style used elsewhere in the sheet
$ws.Range("C23").Copy()
$ws.Range("D26").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to D26 (matches where editing finished)
$ws.Range("D26").Select()
